$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Swap the identity_schema_id / apptyp_code values held in column B.
$ws.Range("B2").Value = 1002
$ws.Range("B3").Value = 1001

# These two cells used a one-off style (numberformat + top alignment);
# normalize them to the same style already used elsewhere in column A
# (plain style, no special number format) so the stray style is dropped.
$ws.Range("A2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("B3").PasteSpecial(-4122)
